# Apply the crypto price/volume update described by the commit diff.
# Values in column D that look like plain numbers are forced back to text
# (matching the source sheet, which stores Price as text) by using the
# classic leading-apostrophe trick, exactly as typing '580.38 in Excel would.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.796.96"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "2.571.25"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'580.38"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "'143.62"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("E9").Value = "  -2.68%  "
$ws.Range("D10").Value = "'5.58"
$ws.Range("E10").Value = "  -1.29%  "
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("E12").Value = "  -2.74%  "
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").Value = "3.031.62"
$ws.Range("D15").Value = "62.709.97"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").Value = "  -3.54%  "
$ws.Range("D17").Value = "2.578.24"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "'11.06"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("D19").Value = "'339.73"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'67.40"
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("E24").Value = "  +5.40%  "
$ws.Range("D25").Value = "'1.59"
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").Value = "'0.164"
$ws.Range("E26").Value = "  -3.69%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'7.98"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").Value = "'8.21"
$ws.Range("E29").Value = "  -3.79%  "
$ws.Range("D30").Value = "'1.91"
$ws.Range("E30").Value = "  -4.39%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0796"
$ws.Range("E31").Value = "  -3.89%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'453.28"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("D34").Value = "'176.51"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("D37").Value = "'18.85"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("D41").Value = "'39.90"
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").Value = "'156.52"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("D43").Value = "'3.68"
$ws.Range("E43").Value = "  -3.65%  "
$ws.Range("D44").Value = "'0.632"
$ws.Range("E44").Value = "  +3.00%  "
$ws.Range("D45").Value = "'20.90"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("D49").Value = "'17.92"
$ws.Range("E49").Value = "  -2.75%  "
$ws.Range("D50").Value = "'11.44"
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("E51").Value = "  -3.89%  "
